$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45131, 0.00705,   0.007102,  0.00705,   0.007066, 0),
    @(45132, 0.007065,  0.007098,  0.00706,   0.007091, 0),
    @(45133, 0.007094,  0.0071445, 0.0070842, 0.0071284, 0),
    @(45134, 0.007131,  0.0072045, 0.007078,  0.0071642, 0),
    @(45135, 0.007172,  0.00724,   0.007081,  0.007081, 0),
    @(45138, 0.007083,  0.007106,  0.007011,  0.007027, 0),
    @(45139, 0.007024,  0.00703,   0.006967,  0.007001, 0),
    @(45140, 0.006991,  0.007029,  0.006972,  0.006976, 0),
    @(45141, 0.0069742, 0.0070146, 0.0069514, 0.007005, 0)
)

$startRow = 780
$endRow = $startRow + $data.Length - 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
}

# Copy the date-column formatting (style index 2) from the last existing
# data row down through all the newly appended rows, matching the
# original workbook's per-row style pattern.
$ws.Range("A779").Copy()
$ws.Range("A$startRow`:A$endRow").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0
